$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("welcome")

# Remove the two now-unused dialogue rows ("notinterested" at row 11 and
# "knowit" at row 5). Delete bottom-up so the row indices used below stay
# valid for the deletion that happens first.
$ws.Rows.Item(11).Delete()
$ws.Rows.Item(5).Delete()

# Update the remaining dialogue text (row numbers below are post-deletion).
$ws.Range("B4").Value = "头也不抬，随意"
$ws.Range("B6").Value = "SOLID哎……人类[困]"
$ws.Range("B7").Value = "想推荐陌生人也去看下《海达·高布乐》，你"
$ws.Range("B8").Value = "SOLID嗯，那我继续看书了，祝今天好心情[调皮]`nHOLD5`nTRANSheddacomeagain"
$ws.Range("B9").Value = "想劝陌生人也去看看《海达·高布乐》，你"

# Match the author's final selection state on the sheet.
$ws.Activate()
$null = $ws.Range("B25").Select()
